$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H6").Value = "hg"
$ws.Range("L6").Value = "fgg"
$ws.Range("F11").Value = "df"
$ws.Range("E8").Value = "td"

$ws.Range("E8").Select()
